$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (category labels): slashes removed, spaces -> underscores, leading slash stripped
$ws.Range("A2").Value = "Independent-Not_for_profit"
$ws.Range("A3").Value = "Government-Local_Authority"
$ws.Range("A4").Value = "Independent-Private"
$ws.Range("A5").Value = "Independent-Unknown"
$ws.Range("A6").Value = "Independent-National_Trust"
$ws.Range("A7").Value = "University"
$ws.Range("A8").Value = "Unknown"
$ws.Range("A9").Value = "Government-National"
$ws.Range("A10").Value = "Independent-English_Heritage"
$ws.Range("A11").Value = "Independent-National_Trust_for_Scotland"
$ws.Range("A12").Value = "Independent-Historic_Environment_Scotland"
$ws.Range("A13").Value = "Government-Other"
$ws.Range("A14").Value = "Government-Cadw"

# Column B (frequency counts) updated for 2022 stats
$ws.Range("B2").Value = 1734
$ws.Range("B3").Value = 922
$ws.Range("B4").Value = 751
$ws.Range("B5").Value = 221
$ws.Range("B9").Value = 82
$ws.Range("B10").Value = 53
$ws.Range("B11").Value = 27
$ws.Range("B12").Value = 21
$ws.Range("B13").Value = 10
$ws.Range("B7").Value = 110
$ws.Range("B8").Value = 110

# Column C (percentages) updated for 2022 stats
$ws.Range("C2").Value = 41
$ws.Range("C3").Value = 21.8
$ws.Range("C4").Value = 17.76
$ws.Range("C5").Value = 5.23
$ws.Range("C6").Value = 4.37
$ws.Range("C7").Value = 2.6
$ws.Range("C8").Value = 2.6
$ws.Range("C9").Value = 1.94
$ws.Range("C10").Value = 1.25
$ws.Range("C11").Value = 0.64
$ws.Range("C12").Value = 0.5
$ws.Range("C13").Value = 0.24
